$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename "Evalua" -> "Eval (ranking)"
# ---------------------------------------------------------------
$evalRanking = $wb.Worksheets.Item("Evalua")
$evalRanking.Name = "Eval (ranking)"

# ---------------------------------------------------------------
# 2. Insert a new worksheet right after "Eval (ranking)" and name it
#    "Eval (binary)". This becomes the new sheet3.xml physically and
#    pushes "Legend" down to sheet4.xml, matching the target layout.
# ---------------------------------------------------------------
$evalBinary = $wb.Worksheets.Add($null, $evalRanking)
$evalBinary.Name = "Eval (binary)"

# ---------------------------------------------------------------
# 3. Update the Legend sheet labels first (to match shared-string
#    ordering of the authored workbook).
# ---------------------------------------------------------------
$legend = $wb.Worksheets.Item("Legend")

# widen column A to match target (bestFit width)
$legend.Columns.Item(1).ColumnWidth = 16.33203125

# rename "Evaluation" header to "Evaluation (ranked)"
$legend.Range("A10").Value = "Evaluation (ranked)"

# add the new "Evaluation (binary)" legend section heading
$legend.Range("A13").Value = "Evaluation (binary)"

# ---------------------------------------------------------------
# 4. Populate "Eval (binary)" sheet with the new table.
# ---------------------------------------------------------------
$evalBinary.Range("A1").Value = "Version"
$evalBinary.Range("B1").Value = "Precision"
$evalBinary.Range("C1").Value = "Recall"
$evalBinary.Range("D1").Value = "F1 score"

$evalBinary.Range("A8").Value = "3(avg)"
$evalBinary.Range("A9").Value = "3(min)"
$evalBinary.Range("A10").Value = "3(max)"

$evalBinary.Range("A2").Value = "1(avg)"
$evalBinary.Range("A3").Value = "1(min)"
$evalBinary.Range("A4").Value = "1(max)"

$evalBinary.Range("A5").Value = "2(avg)"
$evalBinary.Range("A6").Value = "2(min)"
$evalBinary.Range("A7").Value = "2(max)"

$evalBinary.Range("B2").Value = 1
$evalBinary.Range("C2").Value = 0.45
$evalBinary.Range("D2").Value = 0.62

$evalBinary.Range("B3").Value = 1
$evalBinary.Range("C3").Value = 0.01
$evalBinary.Range("D3").Value = 0.02

$evalBinary.Range("B4").Value = 1
$evalBinary.Range("C4").Value = 1
$evalBinary.Range("D4").Value = 1

$evalBinary.Range("B5").Value = 0.37
$evalBinary.Range("C5").Value = 0.58
$evalBinary.Range("D5").Value = 0.45

$evalBinary.Range("B6").Value = 0.62
$evalBinary.Range("C6").Value = 0.02
$evalBinary.Range("D6").Value = 0.03

$evalBinary.Range("B7").Value = 0.1
$evalBinary.Range("C7").Value = 1
$evalBinary.Range("D7").Value = 0.18

$evalBinary.Range("B8").Value = 0.44
$evalBinary.Range("C8").Value = 0.82
$evalBinary.Range("D8").Value = 0.57

$evalBinary.Range("B9").Value = 0.17
$evalBinary.Range("C9").Value = 0.01
$evalBinary.Range("D9").Value = 0.02

$evalBinary.Range("B10").Value = 0.23
$evalBinary.Range("C10").Value = 1
$evalBinary.Range("D10").Value = 0.38

$evalBinary.Range("H8").Select() | Out-Null

# ---------------------------------------------------------------
# 5. Finish the Legend sheet rows describing the binary metrics.
# ---------------------------------------------------------------
$legend.Range("B13").Value = "Precision"
$legend.Range("B14").Value = "Recall"
$legend.Range("B15").Value = "F1 score "

$legend.Range("C13").Value = "proportion of predicted positive cases that are actually positive"
$legend.Range("C14").Value = "proportion of actual positive cases that predictions correctly identify"
$legend.Range("C15").Value = "harmonic mean of precision and recall, balancing both into a single metric"

$legend.Range("G22").Select() | Out-Null

# ---------------------------------------------------------------
# 6. Activate "Eval (binary)" sheet (matches workbookView activeTab).
# ---------------------------------------------------------------
$evalBinary.Activate()
